$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 (I0) and J1 (IF), styled like the existing header row (bold, bordered, centered)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("B1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2-69: column I (I0) and column J (IF)
$data = @{
    2 = @(7,7)
    3 = @(6,6)
    4 = @(6,6)
    5 = @(8,8)
    6 = @(5,5)
    7 = @(5,5)
    8 = @(11,11)
    9 = @(7,7)
    10 = @(6,6)
    11 = @(6,6)
    12 = @(7,7)
    13 = @(6,6)
    14 = @(6,6)
    15 = @(6,6)
    16 = @(8,8)
    17 = @(8,8)
    18 = @(7,7)
    19 = @(6,6)
    20 = @(8,8)
    21 = @(10,10)
    22 = @(5,6)
    23 = @(7,7)
    24 = @(7,7)
    25 = @(7,7)
    26 = @(8,8)
    27 = @(7,7)
    28 = @(7,7)
    29 = @(10,10)
    30 = @(7,7)
    31 = @(7,7)
    32 = @(7,7)
    33 = @(7,7)
    34 = @(7,7)
    35 = @(7,8)
    36 = @(7,7)
    37 = @(4,4)
    38 = @(7,7)
    39 = @(6,7)
    40 = @(7,7)
    41 = @(7,7)
    42 = @(7,7)
    43 = @(9,9)
    44 = @(11,11)
    45 = @(5,5)
    46 = @(6,6)
    47 = @(6,6)
    48 = @(7,7)
    49 = @(6,7)
    50 = @(6,6)
    51 = @(9,9)
    52 = @(7,9)
    53 = @(7,7)
    54 = @(6,7)
    55 = @(5,6)
    56 = @(6,6)
    57 = @(7,8)
    58 = @(5,7)
    59 = @(6,7)
    60 = @(6,7)
    61 = @(7,7)
    62 = @(7,9)
    63 = @(9,9)
    64 = @(6,7)
    65 = @(8,9)
    66 = @(9,9)
    67 = @(6,6)
    68 = @(8,9)
    69 = @(5,6)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
